$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 143, shifting rows 143:239 down to 144:240
$ws.Rows.Item(143).Insert()

# Populate the new row 143 with the new data
$ws.Cells.Item(143, 1).Value = 5
$ws.Cells.Item(143, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(143, 3).Value = "Maule"
$ws.Cells.Item(143, 4).Value = 44824
$ws.Cells.Item(143, 5).Value = 7
$ws.Cells.Item(143, 6).Value = 100112017
$ws.Cells.Item(143, 7).Value = "Apio"
$ws.Cells.Item(143, 8).Value = "Americana (o)"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 300
$ws.Cells.Item(143, 11).Value = 10000
$ws.Cells.Item(143, 12).Value = 10000
$ws.Cells.Item(143, 13).Value = 10000
$ws.Cells.Item(143, 14).Value = "`$/docena de matas"
$ws.Cells.Item(143, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(143, 16).Value = 1667
$ws.Cells.Item(143, 17).Value = 6
$ws.Cells.Item(143, 18).Value = "Hortaliza"
